# ---------------------------------------------------------------------------
# electron-showcase.xlsx update
#  - [json] add new command storeKeys(json,jsonpath,var) to the '#system'
#    lookup sheet (column M), inserted alphabetically before storeValue(...)
#  - [target] remove the obsolete "text" category from the target list
#    (column A) - the whole "text" lookup column is being retired
#  - the "text" lookup column (Y) itself is removed, shifting the
#    subsequent lookup columns (web, webalert, webcookie, ws, ws.async, xml)
#    one column to the left (Z->Y, AA->Z, AB->AA, AC->AB, AD->AC, AE->AD)
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# ---------------------------------------------------------------------------
# 1) column A ("target" list, header on row 1, data A2:A31) - drop "text"
# ---------------------------------------------------------------------------
$targetValues = @()
for ($r = 2; $r -le 31; $r++) {
    $v = $ws.Cells.Item($r, 1).Value()
    if ($v -ne $null -and $v -ne "" -and $v -ne "text") {
        $targetValues += $v
    }
}

for ($i = 0; $i -lt $targetValues.Length; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $targetValues[$i]
}
# the list got one entry shorter (31 -> 30 total incl. header) - clear the
# now-trailing, previously-last row
$ws.Cells.Item(31, 1).ClearContents()

# ---------------------------------------------------------------------------
# 2) column M ("json" list, header on row 1, data M2:M17) - insert
#    storeKeys(json,jsonpath,var) right before storeValue(json,jsonpath,var)
# ---------------------------------------------------------------------------
$jsonValuesOld = @()
for ($r = 2; $r -le 17; $r++) {
    $v = $ws.Cells.Item($r, 13).Value()
    if ($v -ne $null -and $v -ne "") {
        $jsonValuesOld += $v
    }
}

$jsonValues = @()
for ($i = 0; $i -lt $jsonValuesOld.Length; $i++) {
    if ($jsonValuesOld[$i] -eq "storeValue(json,jsonpath,var)") {
        $jsonValues += "storeKeys(json,jsonpath,var)"
    }
    $jsonValues += $jsonValuesOld[$i]
}

for ($i = 0; $i -lt $jsonValues.Length; $i++) {
    $ws.Cells.Item($i + 2, 13).Value = $jsonValues[$i]
}

# ---------------------------------------------------------------------------
# 3) remove the "text" lookup column (Y) entirely and shift the remaining
#    lookup columns (web=Z, webalert=AA, webcookie=AB, ws=AC, ws.async=AD,
#    xml=AE) one column to the left. Column AE ends up unused afterwards.
# ---------------------------------------------------------------------------
$srcCols = @(26, 27, 28, 29, 30, 31)   # Z, AA, AB, AC, AD, AE
$dstCols = @(25, 26, 27, 28, 29, 30)   # Y, Z,  AA,  AB,  AC,  AD

for ($r = 1; $r -le 129; $r++) {
    for ($i = 0; $i -lt $srcCols.Length; $i++) {
        $srcVal = $ws.Cells.Item($r, $srcCols[$i]).Value()
        if ($srcVal -ne $null -and $srcVal -ne "") {
            $ws.Cells.Item($r, $dstCols[$i]).Value = $srcVal
        } else {
            $ws.Cells.Item($r, $dstCols[$i]).ClearContents()
        }
    }
    # the old last column (AE, 31) is now vacated
    $ws.Cells.Item($r, 31).ClearContents()
}

# ---------------------------------------------------------------------------
# 4) update the named ranges so they line up with the new layout
# ---------------------------------------------------------------------------
$wb.Names.Item("target").RefersTo = "='#system'!`$A`$2:`$A`$30"
$wb.Names.Item("json").RefersTo = "='#system'!`$M`$2:`$M`$18"
$wb.Names.Item("web").RefersTo = "='#system'!`$Y`$2:`$Y`$129"
$wb.Names.Item("webalert").RefersTo = "='#system'!`$Z`$2:`$Z`$8"
$wb.Names.Item("webcookie").RefersTo = "='#system'!`$AA`$2:`$AA`$8"
$wb.Names.Item("ws").RefersTo = "='#system'!`$AB`$2:`$AB`$17"
$wb.Names.Item("ws.async").RefersTo = "='#system'!`$AC`$2:`$AC`$8"
$wb.Names.Item("xml").RefersTo = "='#system'!`$AD`$2:`$AD`$27"
